# A4 for jonathan. Samlet script til figurer.
$wb = $excel.ActiveWorkbook

$wsY   = $wb.Worksheets.Item("Y")
$wsQ2P = $wb.Worksheets.Item("Q2P")
$wsM   = $wb.Worksheets.Item("M")

# --- Sheet "Y": update existing values, insert two new rows --------------
$wsY.Range("C2").Value = 0.49
$wsY.Range("C3").Value = 0.49

# Old rows 4-5 (gas/EH/0.5 and Y/EH/0.5) move down to rows 6-7.
$wsY.Range("A6").Value = "gas"
$wsY.Range("B6").Value = "EH"
$wsY.Range("C6").Value = 0.5

$wsY.Range("A7").Value = "Y"
$wsY.Range("B7").Value = "EH"
$wsY.Range("C7").Value = 0.5

# Rows 4-5 become the two new entries.
$wsY.Range("A4").Value = "Y_gas"
$wsY.Range("B4").Value = "Y"
$wsY.Range("C4").Value = 0.01

$wsY.Range("A5").Value = "Y_oil"
$wsY.Range("B5").Value = "Y"
$wsY.Range("C5").Value = 0.01

# --- Sheet "Q2P": append two new rows -------------------------------------
$wsQ2P.Range("A4").Value = "Y_oil"
$wsQ2P.Range("B4").Value = "oil"

$wsQ2P.Range("A5").Value = "Y_gas"
$wsQ2P.Range("B5").Value = "gas"

# --- Sheet "M": update a value --------------------------------------------
$wsM.Range("E3").Value = 2.12

# --- Selections on each sheet (M_sets is left untouched) -----------------
[void]$wsY.Range("E7").Select()
[void]$wsQ2P.Range("B6").Select()

# "M" becomes the active sheet/tab; select its cell last so it sticks.
[void]$wsM.Activate()
[void]$wsM.Range("F14").Select()
